$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, 0, 3, 3),
    @(7, 0, 6, 2),
    @(5, 0, 3, 2),
    @(2, 1, 3, 2),
    @(3, 2, 7, 0),
    @(4, 2, 5, 1),
    @(3, 3, 3, 0),
    @(4, 0, 4, 2),
    @(5, 0, 3, 2),
    @(6, 2, 7, 1),
    @(5, 1, 4, 2),
    @(5, 2, 5, 0),
    @(4, 2, 6, 1),
    @(6, 0, 5, 2),
    @(2, 0, 3, 3),
    @(4, 0, 5, 3),
    @(5, 2, 5, 1),
    @(3, 2, 3, 1),
    @(2, 2, 3, 0),
    @(6, 0, 6, 2),
    @(3, 0, 4, 3),
    @(6, 3, 5, 0),
    @(4, 2, 4, 1),
    @(5, 0, 5, 2),
    @(3, 2, 5, 0)
)

$startRow = 2201
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}

$lastRow = $startRow + $data.Count - 1
$excel.ActiveWindow.ScrollRow = 2200
$ws.Range("A" + ($lastRow + 1)).Select()
